$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet tracks a rolling window of the most recent IPO book-building
# entries (newest at row 2, oldest at row 21). This push adds one new
# entry at the top ("에이치이엠파마(구.에이치이엠)"), which shifts every
# existing record down by one row; the oldest record (previous row 21,
# "미래에셋비전스팩6호") rolls off the bottom of the table.

# Shift existing data rows 2-20 down to rows 3-21 (row 21's old content
# is overwritten / dropped, matching the fixed 20-row window).
$ws.Range("A2:F20").Copy($ws.Range("A3"))

# Write the new entry into the now-vacated row 2.
$ws.Range("A2").Value2 = "에이치이엠파마(구.에이치이엠)"
$ws.Range("B2").Value2 = "2024.08.26~08.30"
$ws.Range("C2").Value2 = "18,000~21,000"
$ws.Range("D2").Value2 = "-"
$ws.Range("E2").Value2 = 12546
$ws.Range("F2").Value2 = "신한투자증권"
